$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + '28.133.59'
$ws.Range("E2").Value = "'" + '  +2.16%  '

# Row 3
$ws.Range("D3").Value = "'" + '1.867.15'
$ws.Range("E3").Value = "'" + '  +1.24%  '

# Row 4
$ws.Range("D4").Value = "'" + '0.9999'
$ws.Range("E4").Value = "'" + '  -0.45%  '

# Row 5
$ws.Range("D5").Value = "'" + '336.79'
$ws.Range("E5").Value = "'" + '  +0.97%  '

# Row 6
$ws.Range("D6").Value = "'" + '0.9994'
$ws.Range("E6").Value = "'" + '  -0.46%  '

# Row 7
$ws.Range("D7").Value = "'" + '0.4702'
$ws.Range("E7").Value = "'" + '  +1.41%  '

# Row 8
$ws.Range("E8").Value = "'" + '  +1.39%  '

# Row 9
$ws.Range("D9").Value = "'" + '46.78'
$ws.Range("E9").Value = "'" + '  +2.16%  '

# Row 10
$ws.Range("D10").Value = "'" + '0.07971'
$ws.Range("E10").Value = "'" + '  +1.16%  '

# Row 11
$ws.Range("D11").Value = "'" + '0.9902'
$ws.Range("E11").Value = "'" + '  -0.59%  '

# Row 12
$ws.Range("D12").Value = "'" + '21.62'
$ws.Range("E12").Value = "'" + '  +1.02%  '

# Row 13
$ws.Range("D13").Value = "'" + '5.973'
$ws.Range("E13").Value = "'" + '  +0.46%  '

# Row 14
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'" + '7.245'
$ws.Range("E14").Value = "'" + '  +1.68%  '

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = "'" + '1.849.60'
$ws.Range("E15").Value = "'" + '  +0.20%  '

# Row 16
$ws.Range("D16").Value = "'" + '91.46'
$ws.Range("E16").Value = "'" + '  +3.35%  '

# Row 17
$ws.Range("D17").Value = "'" + '1.002'
$ws.Range("E17").Value = "'" + '  -0.42%  '

# Row 18
$ws.Range("D18").Value = "'" + '0.00001044'
$ws.Range("E18").Value = "'" + '  +0.90%  '

# Row 19
$ws.Range("D19").Value = "'" + '0.06619'
$ws.Range("E19").Value = "'" + '  -0.37%  '

# Row 20
$ws.Range("D20").Value = "'" + '17.62'
$ws.Range("E20").Value = "'" + '  +2.94%  '

# Row 21
$ws.Range("D21").Value = "'" + '0.9993'
$ws.Range("E21").Value = "'" + '  -0.48%  '

# Row 22
$ws.Range("D22").Value = "'" + '28.121.87'
$ws.Range("E22").Value = "'" + '  +2.12%  '

# Row 23
$ws.Range("D23").Value = "'" + '5.418'
$ws.Range("E23").Value = "'" + '  +0.76%  '

# Row 24
$ws.Range("E24").Value = "'" + '  +1.00%  '

# Row 25
$ws.Range("D25").Value = "'" + '2.287'
$ws.Range("E25").Value = "'" + '  -0.66%  '

# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = "'" + '159.10'
$ws.Range("E26").Value = "'" + '  +0.35%  '

# Row 27
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = "'" + '2.056.02'
$ws.Range("E27").Value = "'" + '  -0.67%  '

# Row 28
$ws.Range("D28").Value = "'" + '19.64'
$ws.Range("E28").Value = "'" + '  +0.90%  '

# Row 29
$ws.Range("D29").Value = "'" + '2.125'
$ws.Range("E29").Value = "'" + '  +0.70%  '

# Row 30
$ws.Range("D30").Value = "'" + '5.492'
$ws.Range("E30").Value = "'" + '  +1.66%  '

# Row 31
$ws.Range("D31").Value = "'" + '119.63'
$ws.Range("E31").Value = "'" + '  -0.07%  '

# Row 32
$ws.Range("D32").Value = "'" + '0.9727'
$ws.Range("E32").Value = "'" + '  -0.12%  '

# Row 33
$ws.Range("D33").Value = "'" + '0.09497'
$ws.Range("E33").Value = "'" + '  +1.04%  '

# Row 34
$ws.Range("D34").Value = "'" + '3.573'
$ws.Range("E34").Value = "'" + '  -0.40%  '

# Row 35
$ws.Range("D35").Value = "'" + '5.328'
$ws.Range("E35").Value = "'" + '  +0.69%  '

# Row 36
$ws.Range("D36").Value = "'" + '1.362'
$ws.Range("E36").Value = "'" + '  +1.73%  '

# Row 37
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'" + '0.06100'
$ws.Range("E37").Value = "'" + '  +1.13%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'" + '0.02258'
$ws.Range("E38").Value = "'" + '  +1.38%  '

# Row 39
$ws.Range("D39").Value = "'" + '8.395'
$ws.Range("E39").Value = "'" + '  +1.26%  '

# Row 40
$ws.Range("D40").Value = "'" + '1.171'
$ws.Range("E40").Value = "'" + '  -0.82%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = "'" + '0.5957'
$ws.Range("E41").Value = "'" + '  +1.25%  '

# Row 42
$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = "'" + '0.9995'
$ws.Range("E42").Value = "'" + '  -0.39%  '

# Row 43
$ws.Range("D43").Value = "'" + '0.1876'
$ws.Range("E43").Value = "'" + '  +0.65%  '

# Row 44
$ws.Range("D44").Value = "'" + '10.28'
$ws.Range("E44").Value = "'" + '  -0.44%  '

# Row 45
$ws.Range("D45").Value = "'" + '1.279'
$ws.Range("E45").Value = "'" + '  +3.67%  '

# Row 46
$ws.Range("D46").Value = "'" + '0.5597'
$ws.Range("E46").Value = "'" + '  +0.21%  '

# Row 47
$ws.Range("D47").Value = "'" + '12.16'
$ws.Range("E47").Value = "'" + '  +0.41%  '

# Row 48
$ws.Range("D48").Value = "'" + '1.974'
$ws.Range("E48").Value = "'" + '  +3.54%  '

# Row 49
$ws.Range("D49").Value = "'" + '0.06880'
$ws.Range("E49").Value = "'" + '  +2.79%  '

# Row 50
$ws.Range("D50").Value = "'" + '111.85'
$ws.Range("E50").Value = "'" + '  +0.85%  '

# Row 51
$ws.Range("D51").Value = "'" + '1.000'
$ws.Range("E51").Value = "'" + '  -32.78%  '
